$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenario")

# --- Fix existing note text (row 6 / SCD0172): add missing comma after "Search" ---
$ws.Range("B6").Value = "1. Sebelum Search, Refresh Halaman Terlebih Dahulu`n2. Data Saat Proses Dedicated Tidak Dapat Di Search"

# --- New scenario rows ---
$ws.Range("A7").Value = "SCD0174"
$ws.Range("B7").Value = "1. Buka Dedicated dan Free dalam 1 iterasi harus menambahkan fungsi back"

$ws.Range("A8").Value = "SCD0175"
$ws.Range("B8").Value = "1. Konfirmasi ketika add cart dari 3x menjadi 2x`n2. Setiap buka Dedicated, Free, Kelolaan, dan Prospek harus menambahkan fungsi back dulu"

$ws.Range("A9").Value = "SCD0176"
$ws.Range("B9").Value = "1. Setiap Search di Dedicated, Free, Kelolaan, dan Prospek Store, Refresh Halaman Terlebih Dahulu`n2. Data Saat Proses di Pipeline Tidak Dapat Di Search, harus klik btn refresh terlebih dahulu"

$ws.Range("A10").Value = "SCD0177"
$ws.Range("B10").Value = "1. Setiap Search di Dedicated, Free, Kelolaan, dan Prospek Store, Refresh Halaman Terlebih Dahulu`n2. Data Saat Proses di Pipeline Tidak Dapat Di Search, harus klik btn refresh terlebih dahulu"

# --- Wrap text for the new multi-line note cells ---
$ws.Range("B8:B10").WrapText = $true

# --- Vertical-center alignment for the whole Scenario column (A2:A10) ---
$ws.Range("A2:A10").VerticalAlignment = -4108

# --- Row heights to match the final layout ---
$ws.Rows.Item(8).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 60
$ws.Rows.Item(10).RowHeight = 60

# --- Column width tweak (Note Rombakan column) ---
$ws.Columns.Item(2).ColumnWidth = 68.5

# --- View / selection state ---
$ws.Range("C12").Select()
